$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-29 Wednesday" "2024-05-30 Thursday"
Replace-Text "921×6=" "144×9="
Replace-Text "358×4=" "824×3="
Replace-Text "486×3=" "680×4="
Replace-Text "423×2=" "684×6="
Replace-Text "883×6=" "901×2="
Replace-Text "759×9=" "510×2="
Replace-Text "129×9=" "571×3="
Replace-Text "767×3=" "658×5="
Replace-Text "233×2=" "250×3="
Replace-Text "404×8=" "656×8="
Replace-Text "211×9=" "710×5="
Replace-Text "528×7=" "896×5="
Replace-Text "927×3=" "528×2="
Replace-Text "749×3=" "809×5="
Replace-Text "522×6=" "348×7="
Replace-Text "215×7=" "296×4="
Replace-Text "275×7=" "599×4="
Replace-Text "435×4=" "874×7="
Replace-Text "731×6=" "982×5="
Replace-Text "142×9=" "721×3="
Replace-Text "919×3=" "896×3="
Replace-Text "116×8=" "779×5="
Replace-Text "381×7=" "891×7="
Replace-Text "783×9=" "537×7="
Replace-Text "240×5=" "857×4="
